$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Mega section header row (row 2) - reuse same shared strings as columns A:G
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = "hp"
$ws.Range("J2").Value = "atk"
$ws.Range("K2").Value = "def"
$ws.Range("L2").Value = "speed"
$ws.Range("M2").Value = "spatk"
$ws.Range("N2").Value = "spdef"
$ws.Range("O2").Value = "BST"

# ---------------------------------------------------------------------------
# Row 3 - base stats for the first mega'd mon
# ---------------------------------------------------------------------------
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 50
$ws.Range("L3").Value = 80
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 110
$ws.Range("O3").Formula = "=SUM(I3:N3)"

# ---------------------------------------------------------------------------
# Row 4 - "max "/"mega " label rows
# ---------------------------------------------------------------------------
# Recreate the pre-existing A4:F4 "max " formula as a genuine shared formula
# (it was a set of independent identical formulas before).
$ws.Range("A4:F4").Formula = '="max " & A$2'

$ws.Range("I4").Formula = '="mega " & I$2'
$ws.Range("J4").Formula = '="mega " & J$2'
$ws.Range("K4").Formula = '="mega " & K$2'
$ws.Range("L4").Formula = '="mega " & L$2'
$ws.Range("M4").Formula = '="mega " & M$2'
$ws.Range("N4").Formula = '="mega " & N$2'
$ws.Range("O4").Value = "BST"

# ---------------------------------------------------------------------------
# Row 5 - "max stat" computed row (A:G) + mega'd stats (I:O)
# ---------------------------------------------------------------------------
$ws.Range("C5:F5").Formula = '=ROUNDDOWN(((((((C3*2))*100/100)+5)*110/100)),0)'

$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 72
$ws.Range("K5").Value = 55
$ws.Range("L5").Value = 103
$ws.Range("M5").Value = 110
$ws.Range("N5").Value = 130
$ws.Range("O5").Formula = "=SUM(I5:N5)"

# ---------------------------------------------------------------------------
# Row 6 - "max " label row for the mega section
# ---------------------------------------------------------------------------
$ws.Range("I6:N6").Formula = '="max " & I$2'
$ws.Range("O6").Value = "all +nature"

# ---------------------------------------------------------------------------
# Row 7 - merged header for second mon + mega'd "max stat" computed row
# ---------------------------------------------------------------------------
$ws.Range("I7").Formula = '=(((((I5 * 2)) * 100 / 100) + 100 + 10))'

$ws.Range("J7").Formula = '=ROUNDDOWN(((((((J5*2))*100/100)+5)*110/100)),0)'
$ws.Range("J7").NumberFormat = "0"

$ws.Range("K7:N7").Formula = '=ROUNDDOWN(((((((K5*2))*100/100)+5)*110/100)),0)'
$ws.Range("K7:N7").NumberFormat = "0"
$ws.Range("O7").Value = "0iv/0ev"

# ---------------------------------------------------------------------------
# Row 8 - stat headers (A:G) + "stat mult" row (I:O)
# ---------------------------------------------------------------------------
$ws.Range("I8").Value = "stat mult"
$ws.Range("J8:N8").Formula = '=J5/J3'
$ws.Range("I8:N8").NumberFormat = "0.000"
$ws.Range("O8").Formula = "=O5-O3"

# ---------------------------------------------------------------------------
# Rows 10 and 16 - recreate the other two "max " shared-formula label rows
# ---------------------------------------------------------------------------
$ws.Range("A10:F10").Formula = '="max " & A$2'
$ws.Range("C11:F11").Formula = '=ROUNDDOWN(((((((C9*2))*100/100)+5)*110/100)),0)'

$ws.Range("A16:F16").Formula = '="max " & A$2'
$ws.Range("C17:F17").Formula = '=ROUNDDOWN(((((((C15*2))*100/100)+5)*110/100)),0)'

# ---------------------------------------------------------------------------
# Column widths for the new mega columns (K:N)
# ---------------------------------------------------------------------------
$ws.Columns("K:N").AutoFit()

# ---------------------------------------------------------------------------
# Selection / active cell, matching the author's final cursor position
# ---------------------------------------------------------------------------
$ws.Range("J8:N8").Select()
